# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" sheets to match the latest site scrape.
# (gh-pages output regenerated at commit 456a3b4)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 7059
$ws1.Range("F4").Value  = 73
$ws1.Range("F6").Value  = 159
$ws1.Range("F7").Value  = 7543
$ws1.Range("F8").Value  = 80
$ws1.Range("F13").Value = 420
$ws1.Range("F14").Value = 157
$ws1.Range("F16").Value = 421
$ws1.Range("F17").Value = 55
$ws1.Range("F18").Value = 53
$ws1.Range("F20").Value = 5382
$ws1.Range("F21").Value = 134
$ws1.Range("F22").Value = 191
$ws1.Range("F23").Value = 821
$ws1.Range("F25").Value = 277

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 7059
$ws4.Range("F4").Value  = 73
$ws4.Range("F6").Value  = 159
$ws4.Range("F7").Value  = 7543
$ws4.Range("F8").Value  = 80
$ws4.Range("F13").Value = 420
$ws4.Range("F14").Value = 157
$ws4.Range("F16").Value = 421
$ws4.Range("F17").Value = 55
$ws4.Range("F18").Value = 53
$ws4.Range("F21").Value = 5382
$ws4.Range("F23").Value = 134
$ws4.Range("F24").Value = 191
$ws4.Range("F25").Value = 821
$ws4.Range("F27").Value = 277
